# Modify the query time of IS-Label on Lubm2U

$wb = $excel.ActiveWorkbook

# Worksheet "Lubm2U" holds the IS-Label row being updated, and becomes the active sheet.
$wsLubm2U = $wb.Worksheets.Item("Lubm2U")

# Update the IS-Label (row 2) query-time values B2:G2
$wsLubm2U.Range("B2").Value = 198.491738
$wsLubm2U.Range("C2").Value = 209.02425299999999
$wsLubm2U.Range("D2").Value = 198.53247300000001
$wsLubm2U.Range("E2").Value = 197.81455600000001
$wsLubm2U.Range("F2").Value = 197.05119400000001
$wsLubm2U.Range("G2").Value = 198.380337

# Activate the Lubm2U sheet and update its selection to B2:G2
$wsLubm2U.Activate()
$wsLubm2U.Range("B2:G2").Select()
